$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reverse the order of the period/value rows (16-24): the workbook was
# re-sorted so the most recent period (2201) appears first (row 16) and the
# oldest (2105) appears last (row 24). Columns E (period) and F (value) are
# the ones that actually move; B/C/D/G/H/I/J stay constant across all rows.
$periods = @("2201", "2112", "2111", "2110", "2109", "2108", "2107", "2106", "2105")
$values  = @(30285, 36342, 36342, 36342, 36342, 36342, 36342, 36342, 21805)

for ($i = 0; $i -lt 9; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
